# Apply the "Added new categorical values for preparation medium and storage
# medium" change described in the commit:
#   - preparation_medium sheet: relabel "Neutral Buffered Formalin (NBF)" to
#     "NBF (Neutral Buffered Formalin)", reorder the list, and add 4 new
#     terms (Trumps fixative, DMEM, 2% PFA/2.5% Glutaraldehyde, Biops buffer)
#   - storage_medium sheet: same relabel, plus 2 new terms (Cyro-EM,
#     2% PFA/2.5% Glutaraldehyde)
#   - update the data validation list ranges so they cover the new rows
#   - bump the pav:createdOn metadata timestamp

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. preparation_medium sheet (final, fully-ordered list of 26 entries)
# ---------------------------------------------------------------------
$prepMedium = $wb.Worksheets.Item("preparation_medium")

$prepMediumNames = @(
    "NBF (Neutral Buffered Formalin)",
    "Allprotect tissue reagent (ALL)",
    "CLARITY hydrogel",
    "Trumps fixative",
    "Inflated (OCT)",
    "DMEM",
    "PFA (Paraformaldehyde)",
    "Fixed frozen OCT (Formalin, sucrose protected)",
    "Unknown",
    "Fresh frozen OCT",
    "2% PFA/2.5% Glutaraldehyde",
    "Bouin's",
    "Methanol",
    "PAXgene tissue kit (PXT)",
    "PBS",
    "Ethanol",
    "Inflated (Agarose)",
    "PLP (Periodate-Lysine-Paraformaldehyde)",
    "MACS tissue storage solution",
    "Fresh frozen CMC",
    "Fresh frozen gelatin",
    "RNAlater",
    "Biops buffer",
    "Fixed frozen OCT (Cytofix/Cytoperm)",
    "None",
    "Fixed frozen OCT (PFA, sucrose protected)"
)

$prepMediumUris = @(
    "http://purl.obolibrary.org/obo/OBIB_0000213",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000118",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000134",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000331",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000123",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C185409",
    "http://purl.obolibrary.org/obo/CHEBI_61538",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000116",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C17998",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000126",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000332",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000140",
    "http://purl.obolibrary.org/obo/CHEBI_17790",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C185113",
    "http://purl.obolibrary.org/obo/OBI_0100046",
    "http://purl.obolibrary.org/obo/CHEBI_16236",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000106",
    "http://purl.bioontology.org/ontology/MESH/C046311",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000105",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000130",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000198",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63348",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000330",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000149",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C41132",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000147"
)

for ($i = 0; $i -lt $prepMediumNames.Length; $i++) {
    $row = $i + 1
    $prepMedium.Cells.Item($row, 1).Value = $prepMediumNames[$i]
    $prepMedium.Cells.Item($row, 2).Value = $prepMediumUris[$i]
}

# ---------------------------------------------------------------------
# 2. storage_medium sheet (final, fully-ordered list of 20 entries)
# ---------------------------------------------------------------------
$storageMedium = $wb.Worksheets.Item("storage_medium")

$storageMediumNames = @(
    "PBS",
    "OCT",
    "NBF (Neutral Buffered Formalin)",
    "Ethanol",
    "Allprotect tissue reagent (ALL)",
    "DMSO (no serum)",
    "MACS tissue storage solution",
    "PFA (Paraformaldehyde)",
    "Tris-EDTA",
    "Unknown",
    "Gelatin",
    "DMSO (serum)",
    "RNAlater",
    "Cyro-EM",
    "FFPE (Paraffin embedded)",
    "CMC",
    "None",
    "2% PFA/2.5% Glutaraldehyde",
    "Methanol",
    "PAXgene tissue kit (PXT)"
)

$storageMediumUris = @(
    "http://purl.obolibrary.org/obo/OBI_0100046",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63523",
    "http://purl.obolibrary.org/obo/OBIB_0000213",
    "http://purl.obolibrary.org/obo/CHEBI_16236",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000118",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000115",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000105",
    "http://purl.obolibrary.org/obo/CHEBI_61538",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000135",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C17998",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C65802",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000125",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63348",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000333",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C143028",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C83594",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C41132",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000332",
    "http://purl.obolibrary.org/obo/CHEBI_17790",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C185113"
)

for ($i = 0; $i -lt $storageMediumNames.Length; $i++) {
    $row = $i + 1
    $storageMedium.Cells.Item($row, 1).Value = $storageMediumNames[$i]
    $storageMedium.Cells.Item($row, 2).Value = $storageMediumUris[$i]
}

# ---------------------------------------------------------------------
# 3. Update data validation ranges to cover the new rows
# ---------------------------------------------------------------------
$sampleBlock = $wb.Worksheets.Item("Sample Block")
$sampleBlock.Range("M2:M1001").Validation.Formula1 = "'preparation_medium'!`$A`$1:`$A`$26"
$sampleBlock.Range("Q2:Q1001").Validation.Formula1 = "'storage_medium'!`$A`$1:`$A`$20"

# ---------------------------------------------------------------------
# 4. Bump the pav:createdOn metadata timestamp
# ---------------------------------------------------------------------
$metadata = $wb.Worksheets.Item(".metadata")
$metadata.Range("C2").Value = "2024-03-12T09:42:58-07:00"
